$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy formatting from B1 to C1 (bold, bordered, centered style)
$ws.Range("B1").Copy()
$ws.Range("C1").PasteSpecial(-4122)

# Row 1: add new date header in column C
$ws.Range("C1").Value = "13-01-2023"

$ws.Range("A2").Value = "1810 Renta variable"
$ws.Range("B2").Value = 5057.13
$ws.Range("C2").Value = 5034.36

$ws.Range("A3").Value = "1822 Raices Valores Negociables"
$ws.Range("B3").Value = 11386.1
$ws.Range("C3").Value = 11336.13

$ws.Range("A4").Value = "Adcap IOL Acciones Argentina"
$ws.Range("B4").Value = 1036.47
$ws.Range("C4").Value = 1033.62

$ws.Range("A5").Value = "Allaria Acciones"
$ws.Range("B5").Value = 2685.23
$ws.Range("C5").Value = 2672.24

$ws.Range("A6").Value = "Alpha Acciones"
$ws.Range("B6").Value = 4044.99
$ws.Range("C6").Value = 4057.08

$ws.Range("A7").Value = "Alpha Mega"
$ws.Range("B7").Value = 14787.53
$ws.Range("C7").Value = 14760.18

$ws.Range("A8").Value = "Alpha planeam equil"
$ws.Range("B8").Value = 839.96
$ws.Range("C8").Value = 838.42

$ws.Range("A9").Value = "Alpha renta balan global"
$ws.Range("B9").Value = 26833.05
$ws.Range("C9").Value = 26726.53

$ws.Range("A10").Value = "Argenfunds"
$ws.Range("B10").Value = 409.57
$ws.Range("C10").Value = 407.39

$ws.Range("A11").Value = "Balanz"
$ws.Range("B11").Value = 2807.4
$ws.Range("C11").Value = 850.32

$ws.Range("A12").Value = "Bull Market"
$ws.Range("B12").Value = 0
$ws.Range("C12").Value = 0

$ws.Range("A13").Value = "Consultatio Acciones Argentina"
$ws.Range("B13").Value = 59038.51
$ws.Range("C13").Value = 58979.07

$ws.Range("A14").Value = "Consultatio Renta Variable"
$ws.Range("B14").Value = 18391.88
$ws.Range("C14").Value = 18387.91

$ws.Range("A15").Value = "FBA Calificado"
$ws.Range("B15").Value = 272.72
$ws.Range("C15").Value = 297.76

$ws.Range("A16").Value = "Fima Acciones"
$ws.Range("B16").Value = 5405.32
$ws.Range("C16").Value = 5384.89

$ws.Range("A17").Value = "Fima PB Acciones"
$ws.Range("B17").Value = 1782.22
$ws.Range("C17").Value = 1797.92

$ws.Range("A18").Value = "Goal Acciones Argentinas"
$ws.Range("B18").Value = 49.69
$ws.Range("C18").Value = 48.64

$ws.Range("A19").Value = "Goal acciones plus"
$ws.Range("B19").Value = 270.38
$ws.Range("C19").Value = 269.61

$ws.Range("A20").Value = "HF Acciones Argentinas"
$ws.Range("B20").Value = 4424
$ws.Range("C20").Value = 4422.83

$ws.Range("A21").Value = "HF Acciones Lideres"
$ws.Range("B21").Value = 5595.17
$ws.Range("C21").Value = 5597.41

$ws.Range("A22").Value = "IAM Renta Variable"
$ws.Range("B22").Value = 1315.87
$ws.Range("C22").Value = 1385.43

$ws.Range("A23").Value = "IEB Value"
$ws.Range("B23").Value = 374.56
$ws.Range("C23").Value = 374.78

$ws.Range("A24").Value = "Lombardi"
$ws.Range("B24").Value = 545.15
$ws.Range("C24").Value = 544.42

$ws.Range("A25").Value = "Megainver"
$ws.Range("B25").Value = 1078.23
$ws.Range("C25").Value = 1077.63

$ws.Range("A26").Value = "Pellegrini Acciones"
$ws.Range("B26").Value = 3267.45
$ws.Range("C26").Value = 3265.66

$ws.Range("A27").Value = "Pionero Acciones"
$ws.Range("B27").Value = 2770.12
$ws.Range("C27").Value = 2765.1

$ws.Range("A28").Value = "Premier Renta Variable"
$ws.Range("B28").Value = 516.44
$ws.Range("C28").Value = 518.48

$ws.Range("A29").Value = "Quinquela Acciones"
$ws.Range("B29").Value = 3492.43
$ws.Range("C29").Value = 3514.85

$ws.Range("A30").Value = "Rofex 20 Renta Variable"
$ws.Range("B30").Value = 2590.2
$ws.Range("C30").Value = 2596.22

$ws.Range("A31").Value = "SBS Acciones Argentina"
$ws.Range("B31").Value = 15003.99
$ws.Range("C31").Value = 15015.07

$ws.Range("A32").Value = "Toronto Trust Multimercado"
$ws.Range("B32").Value = 3422.86
$ws.Range("C32").Value = 2302.79

$ws.Range("A33").Value = "avg"
$ws.Range("B33").Value = 6435.31
$ws.Range("C33").Value = 6331.06

$ws.Range("A34").Value = "total"
$ws.Range("B34").Value = 199494.62
$ws.Range("C34").Value = 196262.74
